$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-StyleLike($targetCell, $srcCell) {
    $targetCell.Style = $srcCell.Style
    $targetCell.HorizontalAlignment = $srcCell.HorizontalAlignment
}

# --- Header block (rows 10-12): column A gets the "mtitleStyle" (centered) ---
foreach ($r in 10,11,12) {
    $cell = $ws.Cells.Item($r, 1)
    Copy-StyleLike $cell $ws.Cells.Item(9, 1)
}

# Max column bumped 28 -> 56
$ws.Range("D10").Value = 56
$ws.Range("E10").Value = 56

# Negative marking changed from -3 to -1
$ws.Range("C11").Value = "-1"

# --- New header row 15: G/H "Student Ans" / "Correct Ans" like D15/E15 ---
Copy-StyleLike $ws.Cells.Item(15, 7) $ws.Cells.Item(15, 4)
$ws.Cells.Item(15, 7).Value = "Student Ans"
Copy-StyleLike $ws.Cells.Item(15, 8) $ws.Cells.Item(15, 5)
$ws.Cells.Item(15, 8).Value = "Correct Ans"

# --- New D/E columns for rows 19-40 (second question block), mirroring A/B ---
Copy-StyleLike $ws.Cells.Item(19, 4) $ws.Cells.Item(19, 1)
Copy-StyleLike $ws.Cells.Item(19, 5) $ws.Cells.Item(19, 2)
$ws.Cells.Item(19, 5).Value = "Option A"
Copy-StyleLike $ws.Cells.Item(20, 4) $ws.Cells.Item(20, 1)
Copy-StyleLike $ws.Cells.Item(20, 5) $ws.Cells.Item(20, 2)
$ws.Cells.Item(20, 5).Value = "Option D"
Copy-StyleLike $ws.Cells.Item(21, 4) $ws.Cells.Item(21, 1)
Copy-StyleLike $ws.Cells.Item(21, 5) $ws.Cells.Item(21, 2)
$ws.Cells.Item(21, 5).Value = "Option B"
Copy-StyleLike $ws.Cells.Item(22, 4) $ws.Cells.Item(22, 1)
Copy-StyleLike $ws.Cells.Item(22, 5) $ws.Cells.Item(22, 2)
$ws.Cells.Item(22, 5).Value = "Option C"
Copy-StyleLike $ws.Cells.Item(23, 4) $ws.Cells.Item(23, 1)
Copy-StyleLike $ws.Cells.Item(23, 5) $ws.Cells.Item(23, 2)
$ws.Cells.Item(23, 5).Value = "Option B"
Copy-StyleLike $ws.Cells.Item(24, 4) $ws.Cells.Item(24, 1)
Copy-StyleLike $ws.Cells.Item(24, 5) $ws.Cells.Item(24, 2)
$ws.Cells.Item(24, 5).Value = "Option C"
Copy-StyleLike $ws.Cells.Item(25, 4) $ws.Cells.Item(25, 1)
Copy-StyleLike $ws.Cells.Item(25, 5) $ws.Cells.Item(25, 2)
$ws.Cells.Item(25, 5).Value = "Option D"
Copy-StyleLike $ws.Cells.Item(26, 4) $ws.Cells.Item(26, 1)
Copy-StyleLike $ws.Cells.Item(26, 5) $ws.Cells.Item(26, 2)
$ws.Cells.Item(26, 5).Value = "Option D"
Copy-StyleLike $ws.Cells.Item(27, 4) $ws.Cells.Item(27, 1)
Copy-StyleLike $ws.Cells.Item(27, 5) $ws.Cells.Item(27, 2)
$ws.Cells.Item(27, 5).Value = "Option A"
Copy-StyleLike $ws.Cells.Item(28, 4) $ws.Cells.Item(28, 1)
Copy-StyleLike $ws.Cells.Item(28, 5) $ws.Cells.Item(28, 2)
$ws.Cells.Item(28, 5).Value = "Option A"
Copy-StyleLike $ws.Cells.Item(29, 4) $ws.Cells.Item(29, 1)
Copy-StyleLike $ws.Cells.Item(29, 5) $ws.Cells.Item(29, 2)
$ws.Cells.Item(29, 5).Value = "Option C"
Copy-StyleLike $ws.Cells.Item(30, 4) $ws.Cells.Item(30, 1)
Copy-StyleLike $ws.Cells.Item(30, 5) $ws.Cells.Item(30, 2)
$ws.Cells.Item(30, 5).Value = "Option A"
Copy-StyleLike $ws.Cells.Item(31, 4) $ws.Cells.Item(31, 1)
Copy-StyleLike $ws.Cells.Item(31, 5) $ws.Cells.Item(31, 2)
$ws.Cells.Item(31, 5).Value = "Option D"
Copy-StyleLike $ws.Cells.Item(32, 4) $ws.Cells.Item(32, 1)
Copy-StyleLike $ws.Cells.Item(32, 5) $ws.Cells.Item(32, 2)
$ws.Cells.Item(32, 5).Value = "Option D"
Copy-StyleLike $ws.Cells.Item(33, 4) $ws.Cells.Item(33, 1)
Copy-StyleLike $ws.Cells.Item(33, 5) $ws.Cells.Item(33, 2)
$ws.Cells.Item(33, 5).Value = "Option B"
Copy-StyleLike $ws.Cells.Item(34, 4) $ws.Cells.Item(34, 1)
Copy-StyleLike $ws.Cells.Item(34, 5) $ws.Cells.Item(34, 2)
$ws.Cells.Item(34, 5).Value = "Option D"
Copy-StyleLike $ws.Cells.Item(35, 4) $ws.Cells.Item(35, 1)
Copy-StyleLike $ws.Cells.Item(35, 5) $ws.Cells.Item(35, 2)
$ws.Cells.Item(35, 5).Value = "Option C"
Copy-StyleLike $ws.Cells.Item(36, 4) $ws.Cells.Item(36, 1)
Copy-StyleLike $ws.Cells.Item(36, 5) $ws.Cells.Item(36, 2)
$ws.Cells.Item(36, 5).Value = "Option D"
Copy-StyleLike $ws.Cells.Item(37, 4) $ws.Cells.Item(37, 1)
Copy-StyleLike $ws.Cells.Item(37, 5) $ws.Cells.Item(37, 2)
$ws.Cells.Item(37, 5).Value = "Option B"
Copy-StyleLike $ws.Cells.Item(38, 4) $ws.Cells.Item(38, 1)
Copy-StyleLike $ws.Cells.Item(38, 5) $ws.Cells.Item(38, 2)
$ws.Cells.Item(38, 5).Value = "Option D"
Copy-StyleLike $ws.Cells.Item(39, 4) $ws.Cells.Item(39, 1)
Copy-StyleLike $ws.Cells.Item(39, 5) $ws.Cells.Item(39, 2)
$ws.Cells.Item(39, 5).Value = "Option A"
Copy-StyleLike $ws.Cells.Item(40, 4) $ws.Cells.Item(40, 1)
Copy-StyleLike $ws.Cells.Item(40, 5) $ws.Cells.Item(40, 2)
$ws.Cells.Item(40, 5).Value = "Option A"

# --- New G/H columns for rows 16-21 (third question block), mirroring A/B ---
Copy-StyleLike $ws.Cells.Item(16, 7) $ws.Cells.Item(16, 1)
Copy-StyleLike $ws.Cells.Item(16, 8) $ws.Cells.Item(16, 2)
$ws.Cells.Item(16, 8).Value = "Option A"
Copy-StyleLike $ws.Cells.Item(17, 7) $ws.Cells.Item(17, 1)
Copy-StyleLike $ws.Cells.Item(17, 8) $ws.Cells.Item(17, 2)
$ws.Cells.Item(17, 8).Value = "Option D"
Copy-StyleLike $ws.Cells.Item(18, 7) $ws.Cells.Item(18, 1)
Copy-StyleLike $ws.Cells.Item(18, 8) $ws.Cells.Item(18, 2)
$ws.Cells.Item(18, 8).Value = "Option D"
Copy-StyleLike $ws.Cells.Item(19, 7) $ws.Cells.Item(19, 1)
Copy-StyleLike $ws.Cells.Item(19, 8) $ws.Cells.Item(19, 2)
$ws.Cells.Item(19, 8).Value = "Option A"
Copy-StyleLike $ws.Cells.Item(20, 7) $ws.Cells.Item(20, 1)
Copy-StyleLike $ws.Cells.Item(20, 8) $ws.Cells.Item(20, 2)
$ws.Cells.Item(20, 8).Value = "Option C"
Copy-StyleLike $ws.Cells.Item(21, 7) $ws.Cells.Item(21, 1)
Copy-StyleLike $ws.Cells.Item(21, 8) $ws.Cells.Item(21, 2)
$ws.Cells.Item(21, 8).Value = "Option D"

Write-Host "edit applied"
